$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title / byline / contact
Replace-Text "Contemporary Challenges in Cybersecurity: A Complex Landscape" "The Intriguing Realm of Science: A Journey Through Its Diverse Fields"
Replace-Text "John Smith" "Emily Johnson"
Replace-Text "jsmith@cybersecurity" "johnsonemily00@gmail"

# Intro paragraph
Replace-Text "In the digital era, where technology has become an intricate part of our lives, cybersecurity has emerged as a pressing global concern" "The world of science is an ever-evolving expanse of knowledge, encompassing diverse disciplines that seek to unravel the mysteries of the natural world"
Replace-Text " The rapid advancements and widespread adoption of the internet, cloud computing, and interconnected devices have dramatically expanded the attack surface, creating a fertile ground for malicious actors" " From the intricate dance of atoms to the intricacies of human biology, from the forces that shape the cosmos to the tapestry of historical events, science offers a lens through which we can understand our universe and ourselves"
Replace-Text " The evolving nature of cybersecurity threats in the contemporary landscape demands a comprehensive understanding of the challenges at hand" " In this essay, we will embark on an enlightening journey through the captivating realms of science, exploring its fundamental concepts, groundbreaking discoveries, and profound applications that have shaped human civilization"

# Cybercrime -> Mathematics paragraph
Replace-Text "The rise of cybercrime, driven by sophisticated cybercriminals and organized crime groups, has resulted in increasingly sophisticated and targeted attacks" "In the realm of mathematics, we delve into the language of numbers, exploring patterns, relationships, and abstract concepts that underpin our understanding of the universe"
Replace-Text " The proliferation of ransomware, advanced persistent threats, and zero-day vulnerabilities poses significant risks to businesses, governments, and individuals alike" " From the elegance of geometry to the power of calculus, mathematics provides a framework for describing and predicting natural phenomena, enabling us to construct models, solve complex problems, and unravel the secrets of the cosmos"
Replace-Text " The potential for large-scale data breaches, financial losses, and disruption of critical infrastructure poses grave threats to societal stability and economic prosperity" " Its applications range from engineering and finance to computer science and music, demonstrating the pervasive influence of mathematics in shaping our world"

# Cyberespionage -> Chemistry paragraph
Replace-Text "Cyberespionage, a clandestine form of digital infiltration, poses a formidable challenge to national security and corporate competitiveness" "As we venture into the realm of chemistry, we encounter the fascinating world of elements and compounds, their interactions, and transformations"
Replace-Text " State-sponsored adversaries employ sophisticated techniques to steal sensitive information, intellectual property, and trade secrets" " Chemistry reveals the fundamental principles governing the formation, structure, and properties of matter, enabling us to understand the composition of substances and the changes they undergo"
Replace-Text " The consequences of cyberespionage reach beyond the theft of data, as it can undermine national security, disrupt economic stability, and lead to political instability" " From the vibrant colors of fireworks to the intricate processes of biological systems, chemistry plays a crucial role in various fields such as medicine, agriculture, and materials science, contributing to our ability to heal diseases, produce food, and create innovative materials"

# Summary paragraph
Replace-Text "Contemporary cybersecurity challenges comprise a complex tapestry of cybercrime, cyberespionage, and the rise of advanced persistent threats" "In this essay, we embarked on an enlightening journey through the diverse fields of science, exploring the fundamental concepts, groundbreaking discoveries, and profound applications that have transformed human understanding and shaped our world"
Replace-Text " To effectively address these challenges, collaboration among stakeholders, including governments, industry, academia, and international organizations, is crucial" " From the elegance of mathematics to the intricacies of chemistry, the fascinating realm of life revealed by biology, and the profound insights gained from medicine, arts, government, history, and politics, science has illuminated the mysteries of the universe and empowered us to address global challenges"
Replace-Text " A commitment to responsible and ethical behavior in cyberspace is paramount in fostering a safer and more secure digital environment for all" " Its transformative influence continues to inspire generations of scientists and innovators, propelling humanity toward a future of progress, prosperity, and sustainability"

Write-Output "done"
